$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 51 (4 Sep 2020): worked 4h on responsive CSS for user settings pages
$ws.Range("A51").Value = 44078
$ws.Range("B51").Value = 4
$ws.Range("C51").Value = "CSS: responsiivisuus käyttäjä asetus sivuja"

# Row 52 (5 Sep 2020): worked 2.5h, html done and css mostly done
$ws.Range("A52").Value = 44079
$ws.Range("B52").Value = 2.5
$ws.Range("C52").Value = "html done ja css pääosin"

# Update the selected cell to reflect where the user left off editing
$ws.Range("A56").Select()
